# Adding details to the excel
# Appends a new "SGD" model-results section (rows 74-91) to the bottom of
# Sheet1, mirroring the layout already used for the other classifiers
# (Random Forest / Gradient Boosting / Support Vector / Latent Dirichlet),
# plus the raw console output (confusion matrix + classification report)
# pasted verbatim into column P.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 74: section title "SGD" (merged, bold/centered header style)
# ---------------------------------------------------------------------
$ws.Range("B74:G74").Font.Size = 16
$ws.Range("B74:G74").HorizontalAlignment = -4108
$ws.Range("B74:G74").Merge()
$ws.Range("B74").Value = "SGD"

# Row 75: blank spacer row (keeps the big-font styling of the section)
$ws.Range("B75:G75").Font.Size = 16

# ---------------------------------------------------------------------
# Rows 76-79: confusion matrix table
# ---------------------------------------------------------------------
$ws.Range("B76:G76").Font.Size = 16
$ws.Range("B76").Value = "#"
$ws.Range("D76").Value = "Rush"
$ws.Range("E76").Value = "Unknown"
$ws.Range("F76").Value = "Venku"

$ws.Range("B77:G77").Font.Size = 16
$ws.Range("B77").Value = 0
$ws.Range("C77").Value = "Rush"
$ws.Range("D77").Value = 19
$ws.Range("E77").Value = 6
$ws.Range("F77").Value = 2

$ws.Range("B78:G78").Font.Size = 16
$ws.Range("B78").Value = 1
$ws.Range("C78").Value = "Unknow"
$ws.Range("D78").Value = 0
$ws.Range("E78").Value = 26
$ws.Range("F78").Value = 0

$ws.Range("B79:G79").Font.Size = 16
$ws.Range("B79").Value = 2
$ws.Range("C79").Value = "Venku"
$ws.Range("D79").Value = 0
$ws.Range("E79").Value = 2
$ws.Range("F79").Value = 28

# Rows 80-81: blank spacer rows
$ws.Range("B80:G80").Font.Size = 16
$ws.Range("B81:G81").Font.Size = 16

# ---------------------------------------------------------------------
# Rows 82-88: classification report table (columns B-G)
# ---------------------------------------------------------------------
$ws.Range("B82:G82").Font.Size = 16
$ws.Range("C82").Value = "#"
$ws.Range("D82").Value = "precision"
$ws.Range("E82").Value = "recall"
$ws.Range("F82").Value = "f1-score"
$ws.Range("G82").Value = "support"

$ws.Range("B83:G83").Font.Size = 16
$ws.Range("C83").Value = 0
$ws.Range("D83").Value = 1
$ws.Range("E83").Value = 0.7
$ws.Range("F83").Value = 0.83
$ws.Range("G83").Value = 27

$ws.Range("B84:G84").Font.Size = 16
$ws.Range("C84").Value = 1
$ws.Range("D84").Value = 0.76
$ws.Range("E84").Value = 1
$ws.Range("F84").Value = 0.87
$ws.Range("G84").Value = 26

$ws.Range("B85:G85").Font.Size = 16
$ws.Range("C85").Value = 2
$ws.Range("D85").Value = 0.93
$ws.Range("E85").Value = 0.93
$ws.Range("F85").Value = 0.93
$ws.Range("G85").Value = 30

$ws.Range("B86:G86").Font.Size = 16
$ws.Range("C86").Value = "accuarcy"
$ws.Range("F86").Value = 0.88
$ws.Range("G86").Value = 83

$ws.Range("B87:G87").Font.Size = 16
$ws.Range("C87").Value = "macro avg"
$ws.Range("D87").Value = 0.9
$ws.Range("E87").Value = 0.88
$ws.Range("F87").Value = 0.88
$ws.Range("G87").Value = 83

$ws.Range("B88:G88").Font.Size = 16
$ws.Range("C88").Value = "weighted avg"
$ws.Range("D88").Value = 0.9
$ws.Range("E88").Value = 0.88
$ws.Range("F88").Value = 0.88
$ws.Range("G88").Value = 83

# ---------------------------------------------------------------------
# Column P (rows 82-91): raw console output pasted verbatim -
# confusion matrix followed by the sklearn classification_report text.
# Left in the default (unstyled) font, as it was pasted separately.
# ---------------------------------------------------------------------
$ws.Range("P82").Value = "[[19  6  2]"
$ws.Range("P83").Value = " [ 0 26  0]"
$ws.Range("P84").Value = " [ 0  2 28]]"
$ws.Range("P85").Value = "              precision    recall  f1-score   support"
$ws.Range("P86").Value = "           0       1.00      0.70      0.83        27"
$ws.Range("P87").Value = "           1       0.76      1.00      0.87        26"
$ws.Range("P88").Value = "           2       0.93      0.93      0.93        30"
$ws.Range("P89").Value = "    accuracy                           0.88        83"
$ws.Range("P90").Value = "   macro avg       0.90      0.88      0.88        83"
$ws.Range("P91").Value = "weighted avg       0.90      0.88      0.88        83"

# ---------------------------------------------------------------------
# View state: scroll down to the newly-added section and leave the
# selection where the author finished editing.
# ---------------------------------------------------------------------
$ws.Range("A59").Select()
$excel.ActiveWindow.ScrollRow = 59
$ws.Range("I91").Select()
